$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target dataset (rows 2-21), reflecting corrected filter/sort of candidate statuses
$data = @(
    @(484, "Cognition AI", "Enterprise Account Executive (US)", "Ankita Sharma", "2nd Interview", 45985),
    @(484, "Cognition AI", "Enterprise Account Executive (US)", "Michael Lin", "2nd Interview", 45986),
    @(484, "Cognition AI", "Enterprise Account Executive (US)", "Sanat Valecha", "3rd Interview", 45993),
    @(484, "Cognition AI", "Enterprise Account Executive (US)", "Shruti Devgon", "2nd Interview", 45987),
    @(656, "UnifyApps", "UnifyApps - US Enterprise AE", "Nave Brar", "1st Interview", 45974),
    @(696, "Cognition AI", "Founding EMEA AE", "Rudi Willner", "3rd Interview", 45982),
    @(696, "Cognition AI", "Founding EMEA AE", "Diana Rautner", "1st Interview", 45982),
    @(702, "Cognition AI", "Forward Deployed Engineer / Sales Engineer (UK)", "Joern Messner", "2nd Interview", 45993),
    @(724, "Rox", "ENT AE San Fran", "Payandeh Ekrami", "CV Sent", 45987),
    @(730, "PointFive", "PointFive SE EST", "Paul Plofchan", "1st Interview", 45967),
    @(740, "Axion Ray", "Enterprise Account Executive (East)", "Joshua Easley", "CV Sent", 45987),
    @(740, "Axion Ray", "Enterprise Account Executive (East)", "Paul Plofchan", "CV Sent", 45987),
    @(743, "Conduct AI", "CS1 Conduct AI - Enterprise Account Executive UK (German speaking)", "LEO ALTINTAS", "1st Interview", 45980),
    @(751, "Materialize", "CS3 Materialize - Enterprise Account Executive x3", "Todd Hixson", "4th Interview", 45992),
    @(751, "Materialize", "CS3 Materialize - Enterprise Account Executive x3", "Danielle Messina", "CV Sent", 45987),
    @(751, "Materialize", "CS3 Materialize - Enterprise Account Executive x3", "Paul Plofchan", "3rd Interview", 45993),
    @(751, "Materialize", "CS3 Materialize - Enterprise Account Executive x3", "Tiffany Cole", "3rd Interview", 45993),
    @(755, "Vega Security", "Enterprise AE (Mid-ATL)", "Brandon Nelson", "CV Sent", 45980),
    @(777, "Adaptive6", "SE Director", "Adam Rafael", "CV Sent", 45981),
    @(815, "Conduct AI", "Enterprise Account Executive UK", "Rudi Willner", "1st Interview", 45992)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
